$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3839.5833
$ws.Range("I74").Value = 3620.8333
$ws.Range("K74").Value = 3620.8333
$ws.Range("M74").Value = -2684.8333
$ws.Range("H75").Value = 39385.668
$ws.Range("J75").Value = 39385.668
$ws.Range("L75").Value = 39385.668
$ws.Range("N75").Value = -41257.668
$ws.Range("H76").Value = 9266349
$ws.Range("I76").Value = 7602.136
$ws.Range("J76").Value = 23815808
$ws.Range("K76").Value = 7602.136
$ws.Range("L76").Value = 23815808
$ws.Range("M76").Value = -7287.136
$ws.Range("N76").Value = -23816438
$ws.Range("H77").Value = 3839.5833
$ws.Range("I77").Value = 3620.8333
$ws.Range("K77").Value = 18104.1665
$ws.Range("M77").Value = -13424.1665
$ws.Range("H78").Value = 39385.668
$ws.Range("J78").Value = 39385.668
$ws.Range("L78").Value = 118157.004
$ws.Range("N78").Value = -127517.004
$ws.Range("H79").Value = 9266349
$ws.Range("I79").Value = 7602.136
$ws.Range("J79").Value = 23815808
$ws.Range("K79").Value = 7602.136
$ws.Range("L79").Value = 23815808
$ws.Range("M79").Value = -6510.136
$ws.Range("N79").Value = -23817992
$ws.Range("H80").Value = 472.3
$ws.Range("I80").Value = 367.27274
$ws.Range("J80").Value = 600.6667
$ws.Range("K80").Value = 1101.81822
$ws.Range("L80").Value = 1802.0001
$ws.Range("M80").Value = -103.8182200000001
$ws.Range("N80").Value = -3798.0001
$ws.Range("H82").Value = 614
$ws.Range("I82").Value = 614
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1842
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1436
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 472.3
$ws.Range("I83").Value = 367.27274
$ws.Range("J83").Value = 600.6667
$ws.Range("K83").Value = 3305.45466
$ws.Range("L83").Value = 5406.0003
$ws.Range("M83").Value = 1686.54534
$ws.Range("N83").Value = -15390.0003
$ws.Range("H85").Value = 614
$ws.Range("I85").Value = 614
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1842
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -438
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 829526.25
$ws.Range("I86").Value = 7053.1113
$ws.Range("J86").Value = 1608711.4
$ws.Range("K86").Value = 7053.1113
$ws.Range("L86").Value = 1608711.4
$ws.Range("M86").Value = -5930.1113
$ws.Range("N86").Value = -1610957.4
$ws.Range("H87").Value = 31000
$ws.Range("J87").Value = 31000
$ws.Range("L87").Value = 31000
$ws.Range("N87").Value = -33496
$ws.Range("H88").Value = 30425114
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 45636420
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 45636420
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -45637232
$ws.Range("H89").Value = 829526.25
$ws.Range("I89").Value = 7053.1113
$ws.Range("J89").Value = 1608711.4
$ws.Range("K89").Value = 35265.5565
$ws.Range("L89").Value = 8043557
$ws.Range("M89").Value = -29649.5565
$ws.Range("N89").Value = -8054789
$ws.Range("H90").Value = 31000
$ws.Range("J90").Value = 31000
$ws.Range("L90").Value = 93000
$ws.Range("N90").Value = -105480
$ws.Range("H91").Value = 30425114
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 45636420
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 45636420
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -45639228
$ws.Range("H105").Value = 30635.5
$ws.Range("J105").Value = 30635.5
$ws.Range("L105").Value = 30635.5
$ws.Range("N105").Value = -37623.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5652312.5
$ws.Range("I32").Value = 1983.1666
$ws.Range("K32").Value = 1983.1666
$ws.Range("M32").Value = -1696.1666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 22729026
$ws.Range("I122").Value = 27779600
$ws.Range("K122").Value = 83338800
$ws.Range("M122").Value = -83336350

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 8821.1
$ws.Range("J115").Value = 10293.294
$ws.Range("L115").Value = 30879.882
$ws.Range("N115").Value = -33229.882
$ws.Range("H132").Value = 7661.0625
$ws.Range("I132").Value = 766.5
$ws.Range("J132").Value = 11797.8
$ws.Range("K132").Value = 6898.5
$ws.Range("L132").Value = 106180.2
$ws.Range("M132").Value = -4368.5
$ws.Range("N132").Value = -111240.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4169109.8
$ws.Range("I80").Value = 2523.75
$ws.Range("J80").Value = 25002040
$ws.Range("K80").Value = 2523.75
$ws.Range("L80").Value = 25002040
$ws.Range("M80").Value = -1525.75
$ws.Range("N80").Value = -25004036
$ws.Range("H83").Value = 4169109.8
$ws.Range("I83").Value = 2523.75
$ws.Range("J83").Value = 25002040
$ws.Range("K83").Value = 12618.75
$ws.Range("L83").Value = 125010200
$ws.Range("M83").Value = -7626.75
$ws.Range("N83").Value = -125020184

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 48780
$ws.Range("J123").Value = 48780
$ws.Range("L123").Value = 48780
$ws.Range("N123").Value = -58580

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 2490
$ws.Range("J26").Value = 3980
$ws.Range("L26").Value = 3980
$ws.Range("N26").Value = -4566
$ws.Range("H29").Value = 4740
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 4740
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 4740
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -5320
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
